$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in this sheet store text (inline strings), including values that
# look numeric (e.g. "594.10", "1.00", dotted-thousands prices like
# "67.179.19"). Force each target cell to Text format before writing so
# Excel does not silently coerce/normalize the string into a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.179.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.666.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.660.37"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.24%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.92%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.41"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -6.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000239"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.288.47"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.672.40"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.197.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.54%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.05%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -7.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000137"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -5.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -5.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -6.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.68"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.52"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.807.16"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.107"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.604.45"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.990"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.74"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.46%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.94%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "434.22"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.63"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.23%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.04%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.36"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -8.24%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.747.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.67%  "
